$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Rows 8 and 9 swap: Cardano moves up to row 8, Dogecoin moves down to row 9
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Range('D8') '0.244'
$ws.Range('E8').Value = '  -0.59%  '

$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range('D9') '0.0606'
$ws.Range('E9').Value = '  -0.09%  '

# Price (D) and Volume(1h) (E) updates
# Row 2
$ws.Range('D2').Value = '26.288.85'
$ws.Range('E2').Value = '  +0.36%  '

# Row 3
$ws.Range('D3').Value = '1.600.66'
$ws.Range('E3').Value = '  +0.90%  '

# Row 5
Set-TextValue $ws.Range('D5') '212.62'
$ws.Range('E5').Value = '  +0.52%  '

# Row 6
Set-TextValue $ws.Range('D6') '0.502'
$ws.Range('E6').Value = '  -0.17%  '

# Row 7
Set-TextValue $ws.Range('D7') '1.00'
$ws.Range('E7').Value = '  +0.43%  '

# Row 10
Set-TextValue $ws.Range('D10') '18.87'
$ws.Range('E10').Value = '  -1.96%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.0854'
$ws.Range('E11').Value = '  +0.84%  '

# Row 12
$ws.Range('D12').Value = '1.828.91'
$ws.Range('E12').Value = '  +1.09%  '

# Row 13
$ws.Range('D13').Value = '1.608.71'
$ws.Range('E13').Value = '  +1.34%  '

# Row 14
Set-TextValue $ws.Range('D14') '4.01'
$ws.Range('E14').Value = '  -0.20%  '

# Row 16
Set-TextValue $ws.Range('D16') '63.54'
$ws.Range('E16').Value = '  -0.83%  '

# Row 17
$ws.Range('D17').Value = '26.319.93'
$ws.Range('E17').Value = '  +0.46%  '

# Row 18
Set-TextValue $ws.Range('D18') '225.54'
$ws.Range('E18').Value = '  +6.00%  '

# Row 19
$ws.Range('D19').Value = '0.0₃0723'
$ws.Range('E19').Value = '  -0.27%  '

# Row 20
Set-TextValue $ws.Range('D20') '7.55'
$ws.Range('E20').Value = '  +2.71%  '

# Row 21
Set-TextValue $ws.Range('D21') '1.01'
$ws.Range('E21').Value = '  +0.50%  '

# Row 22
Set-TextValue $ws.Range('D22') '4.30'
$ws.Range('E22').Value = '  +1.34%  '

# Row 23
Set-TextValue $ws.Range('D23') '2.17'
$ws.Range('E23').Value = '  -0.42%  '

# Row 24
Set-TextValue $ws.Range('D24') '8.96'
$ws.Range('E24').Value = '  -0.07%  '

# Row 25
Set-TextValue $ws.Range('D25') '145.30'
$ws.Range('E25').Value = '  +1.17%  '

# Row 27
Set-TextValue $ws.Range('D27') '6.93'
$ws.Range('E27').Value = '  -0.76%  '

# Row 29
Set-TextValue $ws.Range('D29') '15.41'
$ws.Range('E29').Value = '  +1.86%  '

# Row 30
Set-TextValue $ws.Range('D30') '0.0493'
$ws.Range('E30').Value = '  -0.72%  '

# Row 32
Set-TextValue $ws.Range('D32') '3.21'
$ws.Range('E32').Value = '  +0.72%  '

# Row 33
$ws.Range('D33').Value = '1.441.24'
$ws.Range('E33').Value = '  +8.07%  '

# Row 34
Set-TextValue $ws.Range('D34') '2.96'
$ws.Range('E34').Value = '  +0.80%  '

# Row 35
Set-TextValue $ws.Range('D35') '2.44'
$ws.Range('E35').Value = '  -0.04%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.562'
$ws.Range('E37').Value = '  -3.14%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.822'
$ws.Range('E39').Value = '  +0.45%  '

# Row 40
Set-TextValue $ws.Range('D40') '5.82'
$ws.Range('E40').Value = '  +0.42%  '

# Row 43
Set-TextValue $ws.Range('D43') '0.927'
$ws.Range('E43').Value = '  -2.51%  '

# Row 44
$ws.Range('D44').Value = '1.740.74'
$ws.Range('E44').Value = '  +1.13%  '

# Row 46
Set-TextValue $ws.Range('D46') '60.82'
$ws.Range('E46').Value = '  -0.31%  '

# Row 47
Set-TextValue $ws.Range('D47') '87.26'
$ws.Range('E47').Value = '  +1.43%  '

# Row 49
$ws.Range('D49').Value = '0.0₇0985'
$ws.Range('E49').Value = '  -3.69%  '

# Volume(1h)-only updates
$ws.Range('E15').Value = '  -2.61%  '
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('E28').Value = '  +1.43%  '
$ws.Range('E31').Value = '  +1.15%  '
$ws.Range('E36').Value = '  +0.85%  '
$ws.Range('E38').Value = '  -1.17%  '
$ws.Range('E41').Value = '  +0.45%  '
$ws.Range('E42').Value = '  +1.78%  '
$ws.Range('E45').Value = '  -1.11%  '
$ws.Range('E48').Value = '  +0.75%  '
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('E51').Value = '  +0.43%  '
